# Update "Add data for 2022-03-26" - carjacking by month YoY workbook
# Moves the "through" date forward one day (03-17 -> 03-18) and updates
# the March / Total figures for each year column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet (tab) name reflects the new "as of" date.
$ws.Name = "Through 2022-03-18"

# Shared string label for the March row.
$ws.Range("A4").Value = "March (through 03-18)"

# Updated March counts per year (columns B..I = 2015..2022).
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 33
$ws.Range("E4").Value = 36
$ws.Range("F4").Value = 18
$ws.Range("G4").Value = 36
$ws.Range("H4").Value = 49
$ws.Range("I4").Value = 79

# Updated Total counts per year (columns B..I = 2015..2022).
$ws.Range("B5").Value = 57
$ws.Range("C5").Value = 112
$ws.Range("D5").Value = 164
$ws.Range("E5").Value = 173
$ws.Range("F5").Value = 97
$ws.Range("G5").Value = 177
$ws.Range("H5").Value = 391
$ws.Range("I5").Value = 379
